# Applies the "Add files via upload" revision to the three log sheets
# (accept / reject / font_line): refreshed timestamps + run ids, some
# rows removed, some appended, and the active sheet/selection updated.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    # Writes $text as a literal (non-numeric) string cell even when the
    # text looks like a plain integer, without leaving a stray
    # quote-prefix / number-format style behind on the cell.
    $cell.Formula = '=TEXT("' + $text + '","@")'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# Sheet "accept"
# ---------------------------------------------------------------------
$accept = $wb.Worksheets.Item("accept")

$accept.Cells.Item(2,1).Value = "26-05-2023"
$accept.Cells.Item(2,2).Value = "23:22:23"
$accept.Cells.Item(2,3).Value = 3012000

$accept.Rows("3").ClearContents()

$accept.Cells.Item(4,1).Value = "27-05-2023"
$accept.Cells.Item(4,2).Value = "16:06:07"
Set-TextValue $accept.Cells.Item(4,3) "1240603"

$accept.Cells.Item(5,1).Value = "27-05-2023"
$accept.Cells.Item(5,2).Value = "16:07:17"
Set-TextValue $accept.Cells.Item(5,3) "1240612"

$accept.Cells.Item(6,1).Value = "27-05-2023"
$accept.Cells.Item(6,2).Value = "16:07:53"
Set-TextValue $accept.Cells.Item(6,3) "2596012"

$accept.Cells.Item(7,1).Value = "27-05-2023"
$accept.Cells.Item(7,2).Value = "16:09:02"
Set-TextValue $accept.Cells.Item(7,3) "2596007"

# ---------------------------------------------------------------------
# Sheet "reject"
# ---------------------------------------------------------------------
$reject = $wb.Worksheets.Item("reject")

$reject.Cells.Item(2,1).Value = "26-05-2023"
$reject.Cells.Item(2,2).Value = "23:20:17"
$reject.Cells.Item(2,3).Value = 1240666

$reject.Cells.Item(3,1).Value = "27-05-2023"
$reject.Cells.Item(3,2).Value = "16:06:42"
Set-TextValue $reject.Cells.Item(3,3) "3012004"

$reject.Cells.Item(4,1).Value = "27-05-2023"
$reject.Cells.Item(4,2).Value = "16:08:28"
Set-TextValue $reject.Cells.Item(4,3) "1240614"

$reject.Cells.Item(5,1).Value = "27-05-2023"
$reject.Cells.Item(5,2).Value = "16:09:35"
Set-TextValue $reject.Cells.Item(5,3) "1240609"

$reject.Cells.Item(6,1).Value = "27-05-2023"
$reject.Cells.Item(6,2).Value = "16:11:35"
Set-TextValue $reject.Cells.Item(6,3) "3012008"

$reject.Rows("7:13").ClearContents()

# ---------------------------------------------------------------------
# Sheet "font_line"
# ---------------------------------------------------------------------
$fontLine = $wb.Worksheets.Item("font_line")

$fontLine.Cells.Item(2,1).Value = "25-05-2023"
$fontLine.Cells.Item(2,2).Value = "16:02:43"
$fontLine.Cells.Item(2,3).Value = 2596005

$fontLine.Cells.Item(3,1).Value = "25-05-2024"
$fontLine.Cells.Item(3,2).Value = "16:02:44"
$fontLine.Cells.Item(3,3).Value = 3012004

$fontLine.Cells.Item(4,1).Value = "25-05-2025"
$fontLine.Cells.Item(4,2).Value = "16:02:45"
$fontLine.Cells.Item(4,3).Value = 1240614

$fontLine.Cells.Item(5,1).Value = "25-05-2026"
$fontLine.Cells.Item(5,2).Value = "16:02:46"
$fontLine.Cells.Item(5,3).Value = 3012016

$fontLine.Cells.Item(6,1).Value = "25-05-2027"
$fontLine.Cells.Item(6,2).Value = "16:02:47"
$fontLine.Cells.Item(6,3).Value = 1240601

$fontLine.Cells.Item(7,1).Value = "25-05-2028"
$fontLine.Cells.Item(7,2).Value = "16:02:48"
$fontLine.Cells.Item(7,3).Value = 1240609

$fontLine.Cells.Item(8,1).Value = "25-05-2029"
$fontLine.Cells.Item(8,2).Value = "16:02:49"
$fontLine.Cells.Item(8,3).Value = 2596014

$fontLine.Cells.Item(9,1).Value = "25-05-2030"
$fontLine.Cells.Item(9,2).Value = "16:02:50"
$fontLine.Cells.Item(9,3).Value = 3012008

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping (order matters: the sheet
# selected last becomes the active tab, matching activeTab="0").
# ---------------------------------------------------------------------
$fontLine.Range("E6:F6").Select()
$reject.Range("A3").Select()
$accept.Range("A3").Select()
